$d = $word.ActiveDocument

# --- Simple global text replacements (identical in both bill copies, no ambiguity) ---

# Concessioner's name
$d.Content.Find.Execute("Kyla, Jhan", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Smith, Jane", 2)

# Service address
$d.Content.Find.Execute("Tapun, Dalaguete", $true, $false, $false, $false, $false,
                         $true, 1, $false, "456 Market Road", 2)

# Account number
$d.Content.Find.Execute("Account No. : 00016", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Account No. : None", 2)

# Bill number
$d.Content.Find.Execute("Bill No. 00009", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Bill No. 00007", 2)

# Due date
$d.Content.Find.Execute("2025-06-26", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-06-10", 2)

# --- Meter reading row: Date | Previous | Present | Consumed ---
# Previous (200.0 -> 210.0) and Present (210.0 -> 220.0) share a value mid-flight,
# so a blind global Replace-All would corrupt one of them. Target each table cell
# directly by position instead.

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $t.Cell(5, 1).Range.Text = "2025-06-05"
    $t.Cell(5, 2).Range.Text = "210.0"
    $t.Cell(5, 3).Range.Text = "220.0"
}
